$wb = $excel.ActiveWorkbook

# Sheet1 = "Parameters" -> only selection change (B1 -> C2)
$ws1 = $wb.Worksheets.Item("Parameters")
$ws1.Range("C2").Select()

# Sheet2 = "Properties" -> add new "enabled" column (C) with boolean TRUE values
$ws2 = $wb.Worksheets.Item("Properties")
$ws2.Range("C1").Value = "enabled"
$ws2.Range("C2").Value = $true
$ws2.Range("C3").Value = $true
$ws2.Range("C4").Value = $true
$ws2.Range("C5").Value = $true

$ws2.Range("D24").Select()
